# Update the case locations / public exposure sites table.
# The sheet now has 19 data rows (A2:D20), sorted alphabetically by
# Location, with several existing rows edited and several new rows added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Abbotsford', 'Bodriggy Brewing Company  245 Johnston Street, Abbotsford VIC 3067', '28/12/20 2:50pm-5:30pm', 'Case dined at venue'),
    @('Black Rock', 'Smile Buffalo Thai restaurant  305 Beach Road, Black Rock VIC 3193', '27/12/20 7.30pm - 9pm', 'Case dined for dinner'),
    @('Brighton', 'Royal Brighton Yacht Club - outdoor dining  253 Esplanade, Brighton VIC 3186', '29/12/20 12:00pm-2:00pm', 'Case dined in outdoor restaurant'),
    @('Camberwell', 'Tao Dumplings  1 Evans Place, Camberwell VIC 3124', '29/12/20 12:30pm-1:30pm', 'Case ate at restaurant'),
    @('Cheltenham', 'Two Bob Snob, 256 Charman Road', '22/12/2020 1pm - 2pm', 'Case attended Venue'),
    @('Collingwood', 'Stomping Ground Brewing Company, 100 Gipps Street, Collingwood VIC 3066', '28/12/20 6:00pm-7:30pm', 'Case dined at venue'),
    @('Docklands', 'Melbourne Boat Hire - Yarra River Cruise. 45 Newquay Promenade, Docklands VIC 3008', '28/12/20 11.26am-2:00pm', 'Case attended venue'),
    @('Doveton', 'Holy Family Parish Doveton Catholic  100 Power Road, Doveton VIC 3177', '26/12/20 4:00pm-6:00pm', 'Case attended Spanish Service'),
    @('Glen Waverley', 'Village Century City  285-287 Springvale Road, Glen Waverley VIC 3150', '28/12/20 2:45pm-5:30pm', 'Case attended Gold Class screening Wonder Woman 1984'),
    @('Hampton', 'Merrymen Cafe, 2 Small Street, Hampton VIC', '28-12-2020 1:30pm-2:30pm', 'Case ate in store'),
    @('Keysborough', 'Sikh Temple Keysborough, 198-206 Perry Road, Keysborough', '1/01/21 3:00pm-5:00pm', 'Case visited venue'),
    @('McKinnon', 'Hotlocks By Rachael Hairdresser, 260 McKinnon Road, McKinnon VIC 3204', '23/12/20 4:00pm-6:00pm', 'Case had hair cut in store'),
    @('Melbourne', 'Left Bank Melbourne, 1 Southbank Blvd', '25/12/20 12pm - 3pm', 'Case ate in store'),
    @('Melbourne', 'Melbourne Central Lion Hotel, 211 La Trobe Street', '28/12/2020 10:00pm-12.00am', 'Case attended venue'),
    @('Moorabbin', 'Grape and Grain Liquor Cellars, 14/16 Station St', '21/12/20 2pm - 10pm  22/12/20 10am - 6pm  24/12/20 1pm - 10pm  28/12/20 8.05pm - 8.47pm  29/12/20 12pm - 4pm', 'Case''s workplace'),
    @('Mordialloc', 'Woodlands Golf Club - club bar  109 White Street Mordialloc VIC 3195', '23/12/20 12:30pm-1:30pm', 'Case attended club house bar'),
    @('Mordialloc', 'Woodlands Golf Club - club bar  109 White Street Mordialloc VIC 3195', '28/12/20 4:40pm-5:15pm', 'Case attended club house bar'),
    @('Southbank', 'Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank', '23/12/20 8:00pm-11:00pm', 'Case attended restaurant'),
    @('Springvale', 'IKEA Springvale - Cafe and Restaurant, 917 Princes Hwy', '30/12/20 5:30pm - 6.30pm', 'Case dined at cafe')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}

Write-Output "Updated $($data.Length) data rows (A2:D20)."
